# Applies the "math_L-curve" perturbation: activates the
# optimization_parameters sheet (instead of network_weights), renames the
# "Model" parameter to "production_function" and adds a new "L_curve"
# parameter row right under it, and removes the obsolete "Deletion" row
# from the optimization_parameters sheet.

$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")

# --- optimization_parameters sheet content edits -------------------------

# Row 1 used to repeat the "value" header across C1:F1 - clear those out,
# only A1/B1 stay.
$paramsSheet.Range("C1:F1").ClearContents()

# Rename "Model" -> "production_function" (row 8, column A).
$paramsSheet.Range("A8").Value = "production_function"

# Insert a new row for the L_curve parameter right after the
# production_function/Sigmoid row, shifting everything else down by one.
$paramsSheet.Rows.Item(9).Insert()
$paramsSheet.Range("A9").Value = "L_curve"
$paramsSheet.Range("B9").Value = 1
# Match the numeric formatting used by the other "value" column cells.
$paramsSheet.Range("B2").Copy()
$paramsSheet.Range("B9").PasteSpecial(-4122)

# The old "Sheet" row (now at 16) is followed by the "Deletion" row, which
# the author removed entirely (now at row 17, since we inserted a row above).
$paramsSheet.Rows.Item(17).Delete()

# --- sheet selection / active tab ----------------------------------------

# Previously "network_weights" was the selected/active tab; now it's
# "optimization_parameters", with the whole last row selected.
$paramsSheet.Activate()
$paramsSheet.Rows.Item(17).Select()
